# Add a header row to the iris dataset sheet.
#
# The sheet currently holds the raw iris rows starting at row 1 (no
# header). We insert a new row 1, fill it with the column titles
# (sepal_length, sepal_width, petal_length, petal_width, species),
# center them (matching the existing centered-header column style),
# auto-size the columns for the new, wider header text, and move the
# active selection the way the source workbook ended up (K11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row down by one and leave a blank row 1.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "sepal_length"
$ws.Range("B1").Value = "sepal_width"
$ws.Range("C1").Value = "petal_length"
$ws.Range("D1").Value = "petal_width"
$ws.Range("E1").Value = "species"

# Headers are centered, like the rest of column A-D already is.
$ws.Range("A1:E1").HorizontalAlignment = -4108

# Column contents changed (titles are wider than the numbers/labels
# they sit above), so resize the columns to fit.
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

# Match the saved selection/active cell from the edited workbook.
$null = $ws.Range("K11").Select()
